$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the title text in A32 (the URL/hyperlink in B32 stays the same)
$ws.Range("A32").Value = "Request time off work due after surviving violence"

# Reflect the updated selection shown in the saved file
$ws.Range("C32").Select()
